# ---------------------------------------------------------------------------
# BFoCSbQL workbook update:
#  - split "BFoCSbQL-residential" into "BFoCSbQL-urban-residential" (new,
#    copy of the old residential sheet) and "BFoCSbQL-rural-residential"
#    (renamed from the old residential sheet)
#  - refresh "About" sheet narrative text
#  - relabel the header row / remove the stray orange highlight on the
#    "envelope" row of each of the three quality-level sheets
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the residential sheet, rename the copy + the original
# ---------------------------------------------------------------------------
$residential = $wb.Worksheets.Item("BFoCSbQL-residential")
$residential.Copy($residential)                      # new copy lands right before $residential

$urban = $wb.Worksheets.Item("BFoCSbQL-residential (2)")
$rural = $wb.Worksheets.Item("BFoCSbQL-residential")

$urban.Name = "BFoCSbQL-urban-residential"
$rural.Name = "BFoCSbQL-rural-residential"

# ---------------------------------------------------------------------------
# 2. Refresh the "About" sheet narrative (rows 17-32)
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A17").Value = "This variable captures what BAU fraction of new building components qualify"
$about.Range("A18").Value = "for energy efficient rebates."
$about.Range("A19").Value = ""

$about.Range("A20").Value = "Most of these appliances are produced for use within India, not for export, so"
$about.Range("A21").Value = "production trends and sales trends in appliances lead to essentially the same"
$about.Range("A22").Value = "fraction of components by quality level."

$about.Range("A24").Value = 'For appliances, we categorize only 5-star appliances as "rebate-qualifying" and'
$about.Range("A25").Value = '1-star to 4-star appliances as "standard-compliant" (e.g. non-rebate-qualifying).'

$about.Range("A27").Value = 'For lighting, we categorize LEDs as "rebate-qualifying" and other technologies'
$about.Range("A28").Value = 'as "standard-compliant."'

$about.Range("A30").Value = "While the absolute quanitites would change in the rural-residential market, we "
$about.Range("A31").Value = "assume that the proportion (%) of the rebate-qualifying and standard-compliant"
$about.Range("A32").Value = "appliances would be the same as urban-residential."

# ---------------------------------------------------------------------------
# 3. Relabel header row + drop orange highlight on the "envelope" row for
#    each of the three quality-level data sheets
# ---------------------------------------------------------------------------
$dataSheets = @("BFoCSbQL-urban-residential", "BFoCSbQL-rural-residential", "BFoCSbQL-commercial")
foreach ($name in $dataSheets) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A1").Value = "Dimensionless fraction of total components"
    $ws.Range("A1").WrapText = $true
    $ws.Rows(1).RowHeight = 30

    # "envelope" row (row 4) loses its orange highlight fill
    $ws.Range("B4:C4").Interior.ColorIndex = -4142
}

Write-Output "done"
